$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " November 05 2020"
$ws.Range("C3").Value = "Mumbai won by 57 runs"
$ws.Range("D3").Value = "Delhi Capitals"
$ws.Range("E3").Value = "Mumbai Indians"
$ws.Range("F3").Value = "Daniel Sams "

# These look numeric, but the source data stores them as text (t="str"),
# so force text entry with a leading quote-prefix, then strip the
# resulting "quote prefix" cell style back to Normal so the cell keeps
# the default (unstyled) formatting exactly like the rest of the sheet.
$ws.Range("G3").Value = "'0"
$ws.Range("H3").Value = "'2"
$ws.Range("I3").Value = "'0"
$ws.Range("J3").Value = "'0"
$ws.Range("K3").Value = "'0.00"
$ws.Range("G3:K3").Style = "Normal"
